# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (e.g. "64.320.02"); force text format so
# numeric-looking values (e.g. "578.09") are not reinterpreted as numbers.
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D13","D14","D15","D17","D18","D20","D21","D24","D25","D26","D28","D29","D32","D33","D34","D35","D39","D40","D41","D42","D43","D44","D45","D46","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.320.02'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '2.773.17'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '578.09'
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").Value = '160.94'
$ws.Range("E6").Value = '  +1.28%  '

$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").Value = '0.603'
$ws.Range("E8").Value = '  -1.34%  '

$ws.Range("D9").Value = '0.111'
$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("E10").Value = '  -1.37%  '

$ws.Range("E11").Value = '  +4.52%  '

$ws.Range("E12").Value = '  -0.98%  '

$ws.Range("D13").Value = '3.263.82'
$ws.Range("E13").Value = '  +0.18%  '

$ws.Range("D14").Value = '27.32'
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("D15").Value = '63.992.69'
$ws.Range("E15").Value = '  +0.18%  '

$ws.Range("E16").Value = '  -1.90%  '

$ws.Range("D17").Value = '2.788.39'
$ws.Range("E17").Value = '  +0.49%  '

$ws.Range("D18").Value = '12.25'
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").Value = '361.60'
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").Value = '6.71'
$ws.Range("E21").Value = '  -3.02%  '

$ws.Range("E22").Value = '  +0.54%  '

$ws.Range("E23").Value = '  -5.93%  '

$ws.Range("D24").Value = '65.25'
$ws.Range("E24").Value = '  -2.24%  '

$ws.Range("D25").Value = '0.171'
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("D26").Value = '8.66'
$ws.Range("E26").Value = '  -0.24%  '

$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("D28").Value = '0.0₃0922'
$ws.Range("E28").Value = '  -1.59%  '

$ws.Range("D29").Value = '7.41'
$ws.Range("E29").Value = '  +4.24%  '

$ws.Range("E30").Value = '  -1.08%  '

$ws.Range("E31").Value = '  +10.37%  '

$ws.Range("D32").Value = '166.75'
$ws.Range("E32").Value = '  -1.64%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '1.53'
$ws.Range("E33").Value = '  +4.46%  '

$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '5.03'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").Value = '20.26'
$ws.Range("E35").Value = '  -1.66%  '

$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("E38").Value = '  +0.37%  '

$ws.Range("D39").Value = '351.83'
$ws.Range("E39").Value = '  +5.71%  '

$ws.Range("D40").Value = '6.37'
$ws.Range("E40").Value = '  +3.55%  '

$ws.Range("D41").Value = '4.22'
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").Value = '39.18'
$ws.Range("E42").Value = '  -1.59%  '

$ws.Range("D43").Value = '22.78'
$ws.Range("E43").Value = '  +3.64%  '

$ws.Range("D44").Value = '21.73'
$ws.Range("E44").Value = '  -1.93%  '

$ws.Range("D45").Value = '0.0598'
$ws.Range("E45").Value = '  -0.61%  '

$ws.Range("D46").Value = '137.56'
$ws.Range("E46").Value = '  +0.63%  '

$ws.Range("E47").Value = '  -1.84%  '

$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("E49").Value = '  -0.77%  '

$ws.Range("E50").Value = '  -0.45%  '

$ws.Range("D51").Value = '2.140.62'
$ws.Range("E51").Value = '  +0.30%  '
